$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised values for the 2021-01-01 row (B74/C74)
$ws.Range("B74").Value = 209509.8
$ws.Range("C74").Value = -26743.2

# Add the new 2021-04-01 row (row 75).
# Build the date-label text via a formula in a scratch cell so that Excel's
# automatic text-to-date conversion (triggered by direct .Value assignment
# of a "dd-mm-yyyy"-looking string) is avoided, then copy the computed text
# value (not the formula) into A75, leaving it as plain shared-string text.
$ws.Range("Z1").Formula = "=""01-04-2021"""
$ws.Range("Z1").Copy($ws.Range("A75"))
$ws.Range("Z1").Clear()

$ws.Range("B75").Value = 213282.7
$ws.Range("C75").Value = -16700.7
